# Add "Insulin.pmol.L" column (L) to Sheet2: pmol/L = pg/mL (col K) / 5.808 (insulin MW, kDa)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("L1").Value = "Insulin.pmol.L"

$insulinPmolL = @{
    2 = 498.5082644628099
    3 = 540.3379820936639
    4 = 618.39032369146
    5 = 1743.918732782369
    6 = 242.94111570247935
    7 = 1946.387741046832
    8 = 1126.7565426997246
    9 = 1625.2722107438017
    10 = 1746.7355371900828
    11 = 1188.4757231404958
    12 = 1347.4555785123966
    13 = 232.82042011019283
    14 = 129.6462982093664
    15 = 636.9436983471076
    16 = 415.8056129476584
    17 = 1107.0495867768595
    18 = 1219.7011019283748
    19 = 498.2358815426997
    20 = 3111.914600550964
    21 = 490.86208677685954
    22 = 433.6618457300276
    23 = 618.9075413223139
    24 = 854.5204889807163
    25 = 609.0625
    26 = 651.2835743801654
    27 = 2513.8378099173556
    28 = 614.2489669421489
    29 = 189.10416666666666
    30 = 165.96053719008268
    31 = 77.99542011019284
    32 = 177.27358815426996
    33 = 158.79390495867767
    34 = 98.78498622589532
    35 = 182.16769972451792
    36 = 285.2759986225895
    37 = 143.47727272727272
    38 = 224.58006198347107
    39 = 155.90521694214877
    40 = 155.1810261707989
    41 = 203.47417355371903
    42 = 140.51540977961432
    43 = 267.6515151515151
    44 = 160.95206611570248
    45 = 104.80831611570248
    46 = 177.97486225895318
    47 = 71.93805096418733
    48 = 158.0729683195592
    49 = 268.9195936639118
    50 = 153.73014807162534
    51 = 194.60881542699724
}

foreach ($row in $insulinPmolL.Keys) {
    $ws.Cells.Item([int]$row, 12).Value = $insulinPmolL[$row]
}

$ws.Activate()
$ws.Range("N9").Select()
